$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 131089715
$ws.Range("B8").Value = 58043
$ws.Range("E8").Value = 103021
$ws.Range("F8").Value = 'Talltita'
$ws.Range("G8").Value = 'Poecile montanus'
$ws.Range("H8").Value = '(Conrad von Baldenstein, 1827)'
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = '2'
$ws.Range("Q8").Value = 519475
$ws.Range("R8").Value = 7145090
$ws.Range("AC8").Value = ""

# Row 9
$ws.Range("A9").Value = 131089716
$ws.Range("B9").Value = 57881
$ws.Range("E9").Value = 100049
$ws.Range("F9").Value = 'Spillkråka'
$ws.Range("G9").Value = 'Dryocopus martius'
$ws.Range("Q9").Value = 519470
$ws.Range("R9").Value = 7145085
$ws.Range("AC9").Value = 'Färska hack'

# Row 10
$ws.Range("A10").Value = 131089703
$ws.Range("B10").Value = 57884
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = 'Tretåig hackspett'
$ws.Range("G10").Value = 'Picoides tridactylus'
$ws.Range("H10").Value = '(Linnaeus, 1758)'
$ws.Range("I10").Value = ""
$ws.Range("Q10").Value = 519452
$ws.Range("R10").Value = 7145135
$ws.Range("AC10").Value = 'Ringhack äldre'

# Row 12
$ws.Range("A12").Value = 131089687
$ws.Range("Q12").Value = 519498
$ws.Range("R12").Value = 7144669
$ws.Range("AC12").Value = 'Ringhack äldre'

# Row 14
$ws.Range("A14").Value = 131089698
$ws.Range("Q14").Value = 519635
$ws.Range("R14").Value = 7145132
$ws.Range("AC14").Value = 'Ringhack'

# Row 15
$ws.Range("A15").Value = 131089670
$ws.Range("Q15").Value = 519577
$ws.Range("R15").Value = 7144892
$ws.Range("AC15").Value = 'Ringhack färska och äldre'

# Row 16
$ws.Range("A16").Value = 131089709
$ws.Range("Q16").Value = 519421
$ws.Range("R16").Value = 7144974

# Row 36
$ws.Range("A36").Value = 131089726
$ws.Range("B36").Value = 91828
$ws.Range("E36").Value = 5432
$ws.Range("F36").Value = 'Granticka'
$ws.Range("G36").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H36").Value = ""
$ws.Range("Q36").Value = 519447
$ws.Range("R36").Value = 7144857

# Row 37
$ws.Range("A37").Value = 131089695
$ws.Range("B37").Value = 57884
$ws.Range("E37").Value = 100109
$ws.Range("F37").Value = 'Tretåig hackspett'
$ws.Range("G37").Value = 'Picoides tridactylus'
$ws.Range("H37").Value = '(Linnaeus, 1758)'
$ws.Range("Q37").Value = 519643
$ws.Range("R37").Value = 7145008
$ws.Range("AC37").Value = 'Ringhack'

# Row 38
$ws.Range("A38").Value = 131089717
$ws.Range("B38").Value = 80348
$ws.Range("E38").Value = 6458
$ws.Range("F38").Value = 'Lunglav'
$ws.Range("G38").Value = 'Lobaria pulmonaria'
$ws.Range("H38").Value = '(L.) Hoffm.'
$ws.Range("Q38").Value = 519595
$ws.Range("R38").Value = 7144796
$ws.Range("AC38").Value = ""

# Row 40
$ws.Range("A40").Value = 131089702
$ws.Range("M40").Value = 'gammalt bo'
$ws.Range("Q40").Value = 519473
$ws.Range("R40").Value = 7145177
$ws.Range("AC40").Value = 'Bohål ca 3m upp i grantickerötad granhögstubbe'

# Row 41
$ws.Range("A41").Value = 131089708
$ws.Range("M41").Value = ""
$ws.Range("Q41").Value = 519460
$ws.Range("R41").Value = 7145006
$ws.Range("AC41").Value = 'Ringhack äldre'

# Row 57
$ws.Range("A57").Value = 131089707
$ws.Range("Q57").Value = 519463
$ws.Range("R57").Value = 7145006

# Row 58
$ws.Range("A58").Value = 131089689
$ws.Range("Q58").Value = 519609
$ws.Range("R58").Value = 7144815
$ws.Range("AC58").Value = 'Ringhack färska och äldre'

# Row 59
$ws.Range("A59").Value = 131089713
$ws.Range("Q59").Value = 519350
$ws.Range("R59").Value = 7144788

# Row 60
$ws.Range("A60").Value = 131089711
$ws.Range("Q60").Value = 519399
$ws.Range("R60").Value = 7144945
